$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the tab to reflect the new closed-zones snapshot date
$ws.Name = "2023-03-13"

# Add 4 new rows (32-35), matching the formatting already used by row 31
$ws.Range("A31:M31").Copy()
$ws.Range("A32:M35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill them in bottom-up, as the values were appended to the shared-string
# table in that order
$ws.Range("A35").Value = "37"

$ws.Range("A34").Value = "36"
$ws.Range("B34").Value = 8

$ws.Range("A33").Value = "35"
$ws.Range("B33").Value = "7+11"
$ws.Range("K33").Value = "6+10"

$ws.Range("A32").Value = "31"

[void]$ws.Range("G19").Select()
